$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.174.49"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "2.316.33"
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'532.66"
$ws.Range("E5").Value = "  +2.14%  "
$ws.Range("D6").Value = "'132.22"
$ws.Range("E6").Value = "  -3.03%  "
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").Value = "2.337.51"
$ws.Range("E9").Value = "  -1.10%  "
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").Value = "'5.30"
$ws.Range("E12").Value = "  -2.74%  "
$ws.Range("D13").Value = "'0.345"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").Value = "2.732.61"
$ws.Range("E14").Value = "  -1.38%  "
$ws.Range("D15").Value = "'23.45"
$ws.Range("E15").Value = "  -3.44%  "
$ws.Range("D16").Value = "57.175.21"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("E17").Value = "  -2.45%  "
$ws.Range("D18").Value = "2.327.96"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").Value = "'337.66"
$ws.Range("E19").Value = "  +2.44%  "
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("D21").Value = "'6.92"
$ws.Range("E21").Value = "  +2.82%  "
$ws.Range("D22").Value = "'4.16"
$ws.Range("E22").Value = "  -2.24%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "'61.66"
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").Value = "'0.166"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "'8.72"
$ws.Range("E26").Value = "  +4.84%  "
$ws.Range("D27").Value = "'0.984"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("D29").Value = "'170.54"
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("D31").Value = "0.0{0}0724" -f [char]0x2083
$ws.Range("E31").Value = "  -3.13%  "
$ws.Range("E32").Value = "  -3.26%  "
$ws.Range("D33").Value = "'18.53"
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("D35").Value = "'0.991"
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("E36").Value = "  -3.32%  "
$ws.Range("E37").Value = "  -1.38%  "
$ws.Range("D38").Value = "'0.902"
$ws.Range("E38").Value = "  -2.88%  "
$ws.Range("D39").Value = "'1.59"
$ws.Range("E39").Value = "  +0.72%  "
$ws.Range("D40").Value = "'39.18"
$ws.Range("E40").Value = "  +1.58%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'148.68"
$ws.Range("E41").Value = "  -1.61%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").Value = "'0.377"
$ws.Range("E42").Value = "  -1.68%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'5.38"
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("D44").Value = "'3.60"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("D45").Value = "'280.64"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("D48").Value = "'18.69"
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("E49").Value = "  -1.40%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0217"
$ws.Range("E50").Value = "  -2.13%  "
$ws.Range("B51").Value = "Polygon"
$ws.Range("C51").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D51").Value = "'0.382"
$ws.Range("E51").Value = "  +0.07%  "
